# Case_5_92 vm_pu.xlsx: bump bus0 setpoint from 1.05 pu to 1.02 pu (380 kV slack)
# and refresh the resulting power-flow voltage magnitudes for all buses/timesteps.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.029571681496128
$ws.Range("D2").Value = 1.033743402636385
$ws.Range("E2").Value = 1.029389652301257
$ws.Range("F2").Value = 1.040266992574698
$ws.Range("I2").Value = 1.037072132337225
$ws.Range("J2").Value = 1.034717717581279
$ws.Range("K2").Value = 1.036544873570637
$ws.Range("L2").Value = 1.032203704129455
$ws.Range("M2").Value = 1.043049826529365

# row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.030437294224545
$ws.Range("D3").Value = 1.034386866247325
$ws.Range("E3").Value = 1.030122665514712
$ws.Range("F3").Value = 1.041397131050844
$ws.Range("I3").Value = 1.037314751064253
$ws.Range("J3").Value = 1.035224847584768
$ws.Range("K3").Value = 1.036997822719944
$ws.Range("L3").Value = 1.03274505486513
$ws.Range("M3").Value = 1.043989511971094

# row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030997761954389
$ws.Range("D4").Value = 1.034803464705816
$ws.Range("E4").Value = 1.0305976660809
$ws.Range("F4").Value = 1.042129094127922
$ws.Range("I4").Value = 1.037470576502973
$ws.Range("J4").Value = 1.035552712979833
$ws.Range("K4").Value = 1.037290459089609
$ws.Range("L4").Value = 1.033095369938603
$ws.Range("M4").Value = 1.044597657153158

# row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031233467242056
$ws.Range("D5").Value = 1.034978657436549
$ws.Range("E5").Value = 1.030797520547336
$ws.Range("F5").Value = 1.042436975643458
$ws.Range("I5").Value = 1.037535806001302
$ws.Range("J5").Value = 1.035690479237877
$ws.Range("K5").Value = 1.037413374259371
$ws.Range("L5").Value = 1.033242647318055
$ws.Range("M5").Value = 1.044853346207531

# row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031273048127765
$ws.Range("D6").Value = 1.035008076236229
$ws.Range("E6").Value = 1.030831086580238
$ws.Range("F6").Value = 1.042488679922713
$ws.Range("I6").Value = 1.037546741913052
$ws.Range("J6").Value = 1.035713606763405
$ws.Range("K6").Value = 1.037434005828518
$ws.Range("L6").Value = 1.033267376096236
$ws.Range("M6").Value = 1.044896278965147

# row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031000911128781
$ws.Range("D7").Value = 1.034805805426327
$ws.Range("E7").Value = 1.030600335902565
$ws.Range("F7").Value = 1.04213320741101
$ws.Range("I7").Value = 1.037471449201763
$ws.Range("J7").Value = 1.035554554089068
$ws.Range("K7").Value = 1.037292101918234
$ws.Range("L7").Value = 1.033097337848132
$ws.Range("M7").Value = 1.044601073585941

# row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029864144466255
$ws.Range("D8").Value = 1.033960814788718
$ws.Range("E8").Value = 1.029637233400901
$ws.Range("F8").Value = 1.040648785454986
$ws.Range("I8").Value = 1.037154367359588
$ws.Range("J8").Value = 1.03488916229305
$ws.Range("K8").Value = 1.036698042880887
$ws.Range("L8").Value = 1.032386650348236
$ws.Range("M8").Value = 1.043367375386739

# row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027863812154422
$ws.Range("D9").Value = 1.032473690618435
$ws.Range("E9").Value = 1.027945487087538
$ws.Range("F9").Value = 1.038038332583194
$ws.Range("I9").Value = 1.036586731325236
$ws.Range("J9").Value = 1.033714541074001
$ws.Range("K9").Value = 1.035647814744175
$ws.Range("L9").Value = 1.031134565551133
$ws.Range("M9").Value = 1.041194275363845

# row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026532204416528
$ws.Range("D10").Value = 1.031483605670783
$ws.Range("E10").Value = 1.026821347276188
$ws.Range("F10").Value = 1.036301610413001
$ws.Range("I10").Value = 1.036202359795733
$ws.Range("J10").Value = 1.032930092430662
$ws.Range("K10").Value = 1.034945418780788
$ws.Range("L10").Value = 1.030300063861964
$ws.Range("M10").Value = 1.039746131714944

# row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025956079271037
$ws.Range("D11").Value = 1.031055220925228
$ws.Range("E11").Value = 1.026335475590907
$ws.Range("F11").Value = 1.035550443944746
$ws.Range("I11").Value = 1.0360345191285
$ws.Range("J11").Value = 1.03259010441664
$ws.Range("K11").Value = 1.034640752919968
$ws.Range("L11").Value = 1.029938780843875
$ws.Range("M11").Value = 1.039119214879683

# row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025742152395904
$ws.Range("D12").Value = 1.030896150493732
$ws.Range("E12").Value = 1.026155136088923
$ws.Range("F12").Value = 1.035271554585177
$ws.Range("I12").Value = 1.035971965045819
$ws.Range("J12").Value = 1.032463771052387
$ws.Range("K12").Value = 1.034527508666612
$ws.Range("L12").Value = 1.029804594534633
$ws.Range("M12").Value = 1.038886371283173

# row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02578803717695
$ws.Range("D13").Value = 1.030930269323304
$ws.Range("E13").Value = 1.026193813398063
$ws.Range("F13").Value = 1.035331371534178
$ws.Range("I13").Value = 1.035985392629689
$ws.Range("J13").Value = 1.032490872085949
$ws.Range("K13").Value = 1.034551803448914
$ws.Range("L13").Value = 1.029833377471517
$ws.Range("M13").Value = 1.038936316058543

# row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.02593839454082
$ws.Range("D14").Value = 1.031042071070793
$ws.Range("E14").Value = 1.02632056591178
$ws.Range("F14").Value = 1.035527388264642
$ws.Range("I14").Value = 1.036029352681278
$ws.Range("J14").Value = 1.032579662608443
$ws.Range("K14").Value = 1.034631393696273
$ws.Range("L14").Value = 1.029927688745906
$ws.Range("M14").Value = 1.039099967520191

# row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026031044202417
$ws.Range("D15").Value = 1.031110962581483
$ws.Range("E15").Value = 1.026398680233652
$ws.Range("F15").Value = 1.035648177519571
$ws.Range("I15").Value = 1.036056410031402
$ws.Range("J15").Value = 1.032634363220978
$ws.Range("K15").Value = 1.034680421593919
$ws.Range("L15").Value = 1.029985798428628
$ws.Range("M15").Value = 1.039200801410827

# row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.026570449934688
$ws.Range("D16").Value = 1.031512043195344
$ws.Range("E16").Value = 1.026853611816103
$ws.Range("F16").Value = 1.036351480697129
$ws.Range("I16").Value = 1.036213469266792
$ws.Range("J16").Value = 1.03295264972101
$ws.Range("K16").Value = 1.034965627492416
$ws.Range("L16").Value = 1.030324042413225
$ws.Range("M16").Value = 1.039787741077969

# row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026908931434446
$ws.Range("D17").Value = 1.031763719405702
$ws.Range("E17").Value = 1.027139217199128
$ws.Range("F17").Value = 1.036792870794718
$ws.Range("I17").Value = 1.036311612527778
$ws.Range("J17").Value = 1.033152218301413
$ws.Range("K17").Value = 1.035144390026566
$ws.Range("L17").Value = 1.030536231171008
$ws.Range("M17").Value = 1.040155950269002

# row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027106407222543
$ws.Range("D18").Value = 1.031910549491554
$ws.Range("E18").Value = 1.027305891597008
$ws.Range("F18").Value = 1.037050407750401
$ws.Range("I18").Value = 1.036368722208207
$ws.Range("J18").Value = 1.033268592673775
$ws.Range("K18").Value = 1.035248608591192
$ws.Range("L18").Value = 1.03066000317043
$ws.Range("M18").Value = 1.04037073384009

# row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027173749006574
$ws.Range("D19").Value = 1.031960620076255
$ws.Range("E19").Value = 1.027362737737058
$ws.Range("F19").Value = 1.037138235026506
$ws.Range("I19").Value = 1.03638817211297
$ws.Range("J19").Value = 1.033308268109476
$ws.Range("K19").Value = 1.035284135818828
$ws.Range("L19").Value = 1.030702207191249
$ws.Range("M19").Value = 1.040443971745906

# row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.026872610884053
$ws.Range("D20").Value = 1.031736713650632
$ws.Range("E20").Value = 1.027108565579138
$ws.Range("F20").Value = 1.036745505374429
$ws.Range("I20").Value = 1.03630109670894
$ws.Range("J20").Value = 1.033130809640562
$ws.Range("K20").Value = 1.035125215725398
$ws.Range("L20").Value = 1.030513464709156
$ws.Range("M20").Value = 1.040116443520927

# row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025894116058021
$ws.Range("D21").Value = 1.031009146831148
$ws.Range("E21").Value = 1.02628323672108
$ws.Range("F21").Value = 1.035469662706571
$ws.Range("I21").Value = 1.036016413350627
$ws.Range("J21").Value = 1.032553517294054
$ws.Range("K21").Value = 1.034607958492219
$ws.Range("L21").Value = 1.029899916134033
$ws.Range("M21").Value = 1.039051775654698

# row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025279312637348
$ws.Range("D22").Value = 1.030551990919884
$ws.Range("E22").Value = 1.025765100760674
$ws.Range("F22").Value = 1.034668226505262
$ws.Range("I22").Value = 1.035836203507541
$ws.Range("J22").Value = 1.032190281056922
$ws.Range("K22").Value = 1.034282288732933
$ws.Range("L22").Value = 1.029514213829633
$ws.Range("M22").Value = 1.03838249971959

# row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025605191855622
$ws.Range("D23").Value = 1.030794309571312
$ws.Range("E23").Value = 1.026039699927797
$ws.Range("F23").Value = 1.035093013079233
$ws.Range("I23").Value = 1.035931851438255
$ws.Range("J23").Value = 1.032382864721181
$ws.Range("K23").Value = 1.034454974750134
$ws.Range("L23").Value = 1.029718675888796
$ws.Range("M23").Value = 1.038737283775146

# row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.026889022437627
$ws.Range("D24").Value = 1.031748916291379
$ws.Range("E24").Value = 1.027122415462465
$ws.Range("F24").Value = 1.036766907516877
$ws.Range("I24").Value = 1.036305848773918
$ws.Range("J24").Value = 1.033140483387516
$ws.Range("K24").Value = 1.035133879923169
$ws.Range("L24").Value = 1.030523751875208
$ws.Range("M24").Value = 1.040134294878429

# row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028380607751556
$ws.Range("D25").Value = 1.032857919259744
$ws.Range("E25").Value = 1.028382199703994
$ws.Range("F25").Value = 1.038712568059988
$ws.Range("I25").Value = 1.036734529668756
$ws.Range("J25").Value = 1.034018453912147
$ws.Range("K25").Value = 1.035919723463689
$ws.Range("L25").Value = 1.03145822496993
$ws.Range("M25").Value = 1.041755972015163
